$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

$ws.Range("D5").Value = "Bereken hoeveel regen er valt. "
$ws.Range("E5").Value = "['']"
$ws.Range("F5").Value = 0
$ws.Range("L5").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_edit_4_1763134275.png"
